$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.167.29'
$ws.Range("E2").Value = '  +0.11%  '

$ws.Range("D3").Value = '1.827.82'
$ws.Range("E3").Value = '  -0.52%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9983'
$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.16'
$ws.Range("E5").Value = '  -0.55%  '

$ws.Range("E6").Value = '  -0.81%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07377'
$ws.Range("E8").Value = '  -1.77%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2907'
$ws.Range("E9").Value = '  -1.21%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.11'
$ws.Range("E10").Value = '  -0.77%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07680'
$ws.Range("E11").Value = '  -0.33%  '

$ws.Range("D12").Value = '1.829.19'
$ws.Range("E12").Value = '  -0.70%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.950'
$ws.Range("E13").Value = '  -1.43%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6654'
$ws.Range("E14").Value = '  -1.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.21'
$ws.Range("E15").Value = '  -1.12%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008993'
$ws.Range("E16").Value = '  -3.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.843'
$ws.Range("E17").Value = '  -2.12%  '

$ws.Range("D18").Value = '29.137.62'
$ws.Range("E18").Value = '  -0.04%  '

$ws.Range("D19").Value = '2.070.55'
$ws.Range("E19").Value = '  -1.61%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '236.76'
$ws.Range("E20").Value = '  +1.87%  '

$ws.Range("E21").Value = '  -2.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9995'
$ws.Range("E22").Value = '  -0.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.326'
$ws.Range("E23").Value = '  +2.11%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.0000'
$ws.Range("E24").Value = '  -0.11%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.23'
$ws.Range("E25").Value = '  -1.40%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1410'
$ws.Range("E26").Value = '  +0.31%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.488'
$ws.Range("E27").Value = '  -0.69%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.62'
$ws.Range("E28").Value = '  -1.73%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.485'
$ws.Range("E29").Value = '  -0.82%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05940'
$ws.Range("E30").Value = '  +6.83%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.074'
$ws.Range("E31").Value = '  -1.88%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.085'
$ws.Range("E32").Value = '  -2.39%  '

$ws.Range("E33").Value = '  -0.13%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.860'
$ws.Range("E34").Value = '  +0.34%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7298'
$ws.Range("E35").Value = '  -2.64%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.138'
$ws.Range("E36").Value = '  -0.72%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.606'
$ws.Range("E37").Value = '  -2.30%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.842'
$ws.Range("E38").Value = '  +2.63%  '

$ws.Range("D39").Value = '1.220.57'
$ws.Range("E39").Value = '  -1.74%  '

$ws.Range("E40").Value = '  -2.29%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9194'
$ws.Range("E41").Value = '  +1.83%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.273'
$ws.Range("E42").Value = '  -4.99%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.0000'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.85'
$ws.Range("E44").Value = '  -0.56%  '

$ws.Range("D45").Value = '1.979.03'
$ws.Range("E45").Value = '  -0.89%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.76'
$ws.Range("E46").Value = '  -3.09%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5044'
$ws.Range("E47").Value = '  -0.95%  '

$ws.Range("E48").Value = '  -4.24%  '

$ws.Range("E49").Value = '  -1.71%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.111'
$ws.Range("E50").Value = '  -0.23%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1138'
$ws.Range("E51").Value = '  +2.54%  '
